# DanoPizza assignment workbook edit
# --------------------------------------------------------------------
# The "Pizza_Freq" column (column M) on the "Data View" sheet was
# re-coded from the text labels "N" / "Y" to the numeric values 0 / 1
# (0 = infrequent pizza eater / "N", 1 = frequent pizza eater / "Y").
# This is the only substantive data change in the workbook; every other
# shift in shared-string indices elsewhere in the file is a pure side
# effect of "N" and "Y" no longer being used (and therefore dropped)
# from the shared-strings table once every M-column cell becomes a
# plain number.
# --------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data View")

# row -> new numeric Pizza_Freq value (was "N" -> 0, was "Y" -> 1)
$pizzaFreq = @{
    2  = 0;  3  = 0;  4  = 0;  5  = 1;  6  = 0;
    7  = 0;  8  = 0;  9  = 0;  10 = 1;  11 = 1;
    12 = 1;  13 = 1;  14 = 0;  15 = 0;  16 = 1;
    17 = 1;  18 = 0;  19 = 0;  20 = 0;  21 = 1;
    22 = 0;  23 = 0;  24 = 0;  25 = 1;  26 = 0;
    27 = 1;  28 = 1;  29 = 1;  30 = 1;  31 = 1;
    32 = 1;  33 = 0;  34 = 1;  35 = 0;  36 = 1;
    37 = 0;  38 = 1;  39 = 0;  40 = 1;  41 = 0;
    42 = 0;  43 = 0;  44 = 1;  45 = 1;  46 = 1;
    47 = 0;  48 = 0;  49 = 0;  50 = 1;  51 = 0;
}

foreach ($row in $pizzaFreq.Keys) {
    $ws.Cells.Item($row, 13).Value = $pizzaFreq[$row]
}

# Best-effort match of the recorded UI selection: the author's session
# ended with the whole Pizza_Freq column (M) selected.
$ws.Columns("M").Select()
